$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row containing "RM 232" (row 26) and the row containing
# "SC 92" (row 28). Delete from bottom to top so row indices stay valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()
